$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("kdkab"), shifting the existing
# latitude..status_upload columns one to the right.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "kdkab"

# Match column A's width (custom, not auto-fit) -- use the same
# "ColumnWidth" reading Excel reports for column A. The inserted column
# already inherits column A's header cell style (border/bold/alignment)
# from the insert operation itself.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Update the active selection to match the saved view state.
$ws.Range("D11").Select()
